# Fixed bug with institution list
# Refresh the generated "Bongo" item/test/image rows with a new upload batch
# (new external IDs, new GUIDs, new timestamp, new colors/flavors, new image
# files, new scores/results) and re-point the "Biff" item/test rows at the
# new Bongo external IDs. Also nudges a couple of column widths and the
# active selection on the Biff-Items sheet.

$wb = $excel.ActiveWorkbook

$newTimestamp = "generated 2023-10-12 09:14:28"

# ---------------------------------------------------------------------
# Bongo-Items
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Bongo-Items")

$ws1.Range("A2").Value = "Z00100300023-00082"
$ws1.Range("B2").Value = "3BE35F2B-E261-4DD5-BBF2-0FF6FE330BFC"
$ws1.Range("E2").Value = $newTimestamp
$ws1.Range("F2").Value = "Red"
$ws1.Range("G2").Value = "Strawberry"

$ws1.Range("A3").Value = "Z00100300023-00083"
$ws1.Range("B3").Value = "1116BE19-0257-4F43-AA4E-D4ECD52D2E4E"
$ws1.Range("E3").Value = $newTimestamp
$ws1.Range("F3").Value = "Brown"
$ws1.Range("G3").Value = "Vanilla"

$ws1.Range("A4").Value = "Z00100300023-00084"
$ws1.Range("B4").Value = "6252B910-67E0-4D50-8C57-DEFDD18AC09C"
$ws1.Range("E4").Value = $newTimestamp
$ws1.Range("F4").Value = "Violet"
$ws1.Range("G4").Value = "Raspberry"

$ws1.Range("A5").Value = "Z00100300023-00085"
$ws1.Range("B5").Value = "63E9366A-0320-42C7-9C3B-996F655B9298"
$ws1.Range("E5").Value = $newTimestamp
$ws1.Range("F5").Value = "Yellow"
$ws1.Range("G5").Value = "Vanilla"

$ws1.Columns.Item(1).ColumnWidth = 21.333333333333336
$ws1.Columns.Item(2).ColumnWidth = 12.833333333333332

# ---------------------------------------------------------------------
# Bongo-Item-Images
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Bongo-Item-Images")

$ws2.Range("A2").Value = "Z00100300023-00082"
$ws2.Range("B2").Value = "3BE35F2B-E261-4DD5-BBF2-0FF6FE330BFC"
$ws2.Range("C2").Value = $newTimestamp
$ws2.Range("D2").Value = "../images/laughingcat.jpeg"

$ws2.Range("A3").Value = "Z00100300023-00083"
$ws2.Range("B3").Value = "1116BE19-0257-4F43-AA4E-D4ECD52D2E4E"
$ws2.Range("C3").Value = $newTimestamp
$ws2.Range("D3").Value = "../iamges/apple.jpeg"

$ws2.Range("A4").Value = "Z00100300023-00084"
$ws2.Range("B4").Value = "6252B910-67E0-4D50-8C57-DEFDD18AC09C"
$ws2.Range("C4").Value = $newTimestamp
$ws2.Range("D4").Value = "../images/raccoon.jpeg"

$ws2.Range("A5").Value = "Z00100300023-00085"
$ws2.Range("B5").Value = "63E9366A-0320-42C7-9C3B-996F655B9298"
$ws2.Range("C5").Value = $newTimestamp
$ws2.Range("D5").Value = "../images/lightbulb.jpeg"

# ---------------------------------------------------------------------
# Bongo-Tests
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Bongo-Tests")

$ws3.Range("A2").Value = "Z00100300023-00082"
$ws3.Range("B2").Value = "3BE35F2B-E261-4DD5-BBF2-0FF6FE330BFC"
$ws3.Range("C2").Value = $newTimestamp
$ws3.Range("D2").Value = 77
$ws3.Range("E2").Value = "Fail"

$ws3.Range("A3").Value = "Z00100300023-00083"
$ws3.Range("B3").Value = "1116BE19-0257-4F43-AA4E-D4ECD52D2E4E"
$ws3.Range("C3").Value = $newTimestamp
$ws3.Range("D3").Value = 48
$ws3.Range("E3").Value = "Fail"

$ws3.Range("A4").Value = "Z00100300023-00084"
$ws3.Range("B4").Value = "6252B910-67E0-4D50-8C57-DEFDD18AC09C"
$ws3.Range("C4").Value = $newTimestamp
$ws3.Range("D4").Value = 85
$ws3.Range("E4").Value = "Pass"

$ws3.Range("A5").Value = "Z00100300023-00085"
$ws3.Range("B5").Value = "63E9366A-0320-42C7-9C3B-996F655B9298"
$ws3.Range("C5").Value = $newTimestamp
$ws3.Range("D5").Value = 42
$ws3.Range("E5").Value = "Fail"

# ---------------------------------------------------------------------
# Bongo-Test-Images
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Bongo-Test-Images")

$ws4.Range("A2").Value = "Z00100300023-00082"
$ws4.Range("B2").Value = "3BE35F2B-E261-4DD5-BBF2-0FF6FE330BFC"
$ws4.Range("C2").Value = $newTimestamp
$ws4.Range("D2").Value = "../images/broccoli.jpeg"

$ws4.Range("A3").Value = "Z00100300023-00083"
$ws4.Range("B3").Value = "1116BE19-0257-4F43-AA4E-D4ECD52D2E4E"
$ws4.Range("C3").Value = $newTimestamp
$ws4.Range("D3").Value = "../images/laughingcat.jpeg"

$ws4.Range("A4").Value = "Z00100300023-00084"
$ws4.Range("B4").Value = "6252B910-67E0-4D50-8C57-DEFDD18AC09C"
$ws4.Range("C4").Value = $newTimestamp
$ws4.Range("D4").Value = "../images/dice.jpeg"

$ws4.Range("A5").Value = "Z00100300023-00085"
$ws4.Range("B5").Value = "63E9366A-0320-42C7-9C3B-996F655B9298"
$ws4.Range("C5").Value = $newTimestamp
$ws4.Range("D5").Value = "../images/raccoon.jpeg"

# ---------------------------------------------------------------------
# Biff-Items
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Biff-Items")

$ws5.Range("A2").Value = "Z00100300022-00046"
$ws5.Range("B2").Value = "5D330D27-AA7C-46E6-A90F-ABE6980437D7"
$ws5.Range("F2").Value = 81.34
$ws5.Range("G2").Value = 81.64
$ws5.Range("H2").Value = 83.84
$ws5.Range("I2").Value = "Z00100300023-00084"
$ws5.Range("J2").Value = "Z00100300023-00085"

$ws5.Range("A3").Value = "Z00100300022-00047"
$ws5.Range("B3").Value = "30A495F0-3094-4D66-B131-57715A15ECEB"
$ws5.Range("F3").Value = 86.03
$ws5.Range("G3").Value = 106.24
$ws5.Range("H3").Value = 84.53
$ws5.Range("I3").Value = "Z00100300023-00084"
$ws5.Range("J3").Value = "Z00100300023-00085"

$ws5.Columns.Item(9).ColumnWidth = 29.0
$ws5.Columns.Item(10).ColumnWidth = 24.5

$ws5.Range("J3").Select()

# ---------------------------------------------------------------------
# Biff-Tests
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("Biff-Tests")

$ws6.Range("A2").Value = "Z00100300022-00046"
$ws6.Range("B2").Value = "5D330D27-AA7C-46E6-A90F-ABE6980437D7"
$ws6.Range("C2").Value = $newTimestamp
$ws6.Range("D2").Value = 56
$ws6.Range("E2").Value = "NA"
$ws6.Range("F2").Value = "Fail"

$ws6.Range("A3").Value = "Z00100300022-00047"
$ws6.Range("B3").Value = "30A495F0-3094-4D66-B131-57715A15ECEB"
$ws6.Range("C3").Value = $newTimestamp
$ws6.Range("D3").Value = 91
$ws6.Range("E3").Value = "NA"
$ws6.Range("F3").Value = "Pass"
